$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "cl"
$ws.Range("L1").Value = "br"
$ws.Range("M1").Value = "na"
$ws.Range("N1").Value = "ph_lab"
$ws.Range("O1").Value = "ph_field"
$ws.Range("K2").Value = "mg/l"
$ws.Range("L2").Value = "mg/l"
$ws.Range("M2").Value = "mg/l"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = "-"
$ws.Range("K3").Value = 16
$ws.Range("L3").Value = 50
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 5
$ws.Range("O3").Value = 5.3
$ws.Range("K4").Value = 19
$ws.Range("L4").Value = 224
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 7.6
$ws.Range("K5").Value = 31
$ws.Range("L5").Value = 340
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 8.2
$ws.Range("K6").Value = 34
$ws.Range("L6").Value = 160
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 4.8
$ws.Range("O6").Value = 5
$ws.Range("K7").Value = 41
$ws.Range("L7").Value = 160
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1.8
$ws.Range("K8").Value = 41
$ws.Range("L8").Value = 110
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 1.8
$ws.Range("K9").Value = 44
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 1.8
$ws.Range("K10").Value = 50
$ws.Range("L10").Value = 140
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = 1.8
$ws.Range("K11").Value = 45
$ws.Range("L11").Value = 115
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 2
$ws.Range("O11").Value = 1.8
$ws.Range("K12").Value = 33
$ws.Range("L12").Value = 240
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 1.8
$ws.Range("K13").Value = 29
$ws.Range("L13").Value = 200
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 2
$ws.Range("O13").Value = 1.8
$ws.Range("K14").Value = 33
$ws.Range("L14").Value = 740
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 2
$ws.Range("O14").Value = 1.8
$ws.Range("K15").Value = 35
$ws.Range("L15").Value = 1750
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 2
$ws.Range("O15").Value = 1.8
$ws.Range("K16").Value = 60
$ws.Range("L16").Value = 10800
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = 2.4
$ws.Range("K17").Value = 390
$ws.Range("L17").Value = 2556
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 7.8
$ws.Range("K18").Value = 1320
$ws.Range("L18").Value = 212
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 9
$ws.Range("O18").Value = 8.5
$ws.Range("K19").Value = 5285
$ws.Range("L19").Value = 6845
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 9.9

$null = $ws.Range("O20").Select()
